# Swap the data held in rows 2-5 with the data held in rows 6-9.
# (Row 10 and below, and row 1 header, are left untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$blockA = $ws.Range("A2:AY5").Value()
$blockB = $ws.Range("A6:AY9").Value()

# Columns Y (Startdatum) and AA (Slutdatum) hold plain text dates
# ("yyyy-mm-dd") in the source data; force text format on those columns
# so the COM Value setter doesn't silently coerce them into real dates.
$ws.Range("Y2:Y9").NumberFormat = "@"
$ws.Range("AA2:AA9").NumberFormat = "@"

$ws.Range("A2:AY5").Value = $blockB
$ws.Range("A6:AY9").Value = $blockA

# Drop the temporary Text number format again so the cells end up with
# no explicit style, matching the plain inline-string cells before the edit.
$ws.Range("Y2:Y9").ClearFormats()
$ws.Range("AA2:AA9").ClearFormats()
